# Update the sample Invoice No. / PO No. / Supplier Name / Buyer Name values
# on row 2 of both historical templates.
#
# "Historical Invoice Template": B=PO No., C=Invoice No., E=Supplier Name, F=Buyer Name
# "Historical PO Template":      B=Invoice No., C=PO No., E=Supplier Name, F=Buyer Name

$wb = $excel.ActiveWorkbook

$wsInvoice = $wb.Worksheets.Item("Historical Invoice Template")
$wsInvoice.Range("B2").Value = "Invoice1665029"
$wsInvoice.Range("C2").Value = "Invoice1809241"
$wsInvoice.Range("E2").Value = "Frances6b1j"
$wsInvoice.Range("F2").Value = "Henryd777"

$wsPO = $wb.Worksheets.Item("Historical PO Template")
$wsPO.Range("B2").Value = "Invoice1809241"
$wsPO.Range("C2").Value = "Invoice1665029"
$wsPO.Range("E2").Value = "Frances6b1j"
$wsPO.Range("F2").Value = "Henryd777"
